$p = $ppt.ActivePresentation
$s = $p.Slides.Item(16)
$shape = $s.Shapes.Item(3)
$table = $shape.Table
$table.ApplyStyle("{86177154-FD69-402C-82E4-80446896644B}")
